# Remove the Defect ID's (DE...) and their associated comments from the
# "RTM" sheet's "Comments Failed US- Defect ID Not Testable Comment" column (J)
# per the Release Agent's request.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM")

$rows = @(3, 4, 5, 21, 23, 28, 29, 30, 31, 32, 33, 34, 35, 36, 40)

foreach ($r in $rows) {
    $ws.Range("J$r").Value = $null
}
